{"js": "// \"temporary fix for LibreOffice users\"\n//\n// Adds two new character styles (H4, H5 - H5 based on H4, both bold,\n// H5 also italic) to the document and marks the lone empty paragraph's\n// paragraph-mark with the H5 character style so LibreOffice renders the\n// (currently empty) heading placeholder with the right run formatting.\n\n// --- Add the H4 / H5 character styles --------------------------------\ncontext.document.addStyle(\"H4\", \"Character\");\nawait context.sync();\nconst h4 = context.document.getStyles().getByNameOrNullObject(\"H4\");\nh4.baseStyle = \"DefaultParagraphFont\";\nh4.priority = 1;\nh4.quickStyle = true;\nh4.font.bold = true;\nawait context.sync();\n\ncontext.document.addStyle(\"H5\", \"Character\");\nawait context.sync();\nconst h5 = context.document.getStyles().getByNameOrNullObject(\"H5\");\nh5.baseStyle = \"H4\";\nh5.priority = 1;\nh5.quickStyle = true;\nh5.font.bold = true;\nh5.font.italic = true;\nawait context.sync();\n\n// --- Stamp the document's (only, empty) paragraph mark with the H5 style --\n// We want <w:pPr><w:rPr><w:rStyle w:val=\"H5\"/></w:rPr></w:pPr> on the\n// paragraph. The Word JS API has no direct \"apply character style to the\n// paragraph mark of an empty paragraph\" call (Range.style on a whole-\n// paragraph range sets the *paragraph* style instead), so rebuild the\n// paragraph's XML directly (preserving its existing paragraph-mark\n// attributes, e.g. w14:paraId/rsid*) and push it back in with insertOoxml.\nconst body = context.document.body;\nconst wholeRange = body.getRange(\"Whole\");\nconst existingOoxml = wholeRange.getOoxml();\nawait context.sync();\n\nconst full = existingOoxml.value;\nconst partMarker = 'pkg:name=\"/word/document.xml\"';\nconst partIdx = full.indexOf(partMarker);\nconst bodyIdx = full.indexOf(\"<w:body>\", partIdx);\nconst pStart = full.indexOf(\"<w:p\", bodyIdx);\nconst pEnd = full.indexOf(\">\", pStart);\nconst openTag = full.substring(pStart, pEnd + 1);\n\nconst selfClosing = openTag.endsWith(\"/>\");\nconst attrs = selfClosing\n  ? openTag.substring(4, openTag.length - 2).trim()\n  : openTag.substring(4, openTag.length - 1).trim();\n\nconst newParaXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n<w:body>\n<w:p ${attrs}><w:pPr><w:rPr><w:rStyle w:val=\"H5\"/></w:rPr></w:pPr></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nwholeRange.insertOoxml(newParaXml, \"Replace\");\nawait context.sync();\n", "ps1": "# \"temporary fix for LibreOffice users\"\n#\n# Adds two new character styles (H4, H5 - H5 based on H4, both bold,\n# H5 also italic) to the document and marks the lone empty paragraph's\n# paragraph-mark with the H5 character style so LibreOffice renders the\n# (currently empty) heading placeholder with the right run formatting.\n\n$d = $word.ActiveDocument\n\n# --- Add the H4 / H5 character styles -------------------------------------\n$h4 = $d.Styles.Add(\"H4\", 2)          # 2 = wdStyleTypeCharacter\n$h4.BaseStyle  = \"DefaultParagraphFont\"\n$h4.Priority   = 1\n$h4.QuickStyle = $true\n$h4.Font.Bold  = $true\n\n$h5 = $d.Styles.Add(\"H5\", 2)          # 2 = wdStyleTypeCharacter\n$h5.BaseStyle   = \"H4\"\n$h5.Priority    = 1\n$h5.QuickStyle  = $true\n$h5.Font.Bold   = $true\n$h5.Font.Italic = $true\n\n# --- Stamp the document's (only, empty) paragraph mark with the H5 style --\n# We want <w:pPr><w:rPr><w:rStyle w:val=\"H5\"/></w:rPr></w:pPr> on the\n# paragraph. Range.CharacterStyle lands in the right place (w:pPr/w:rPr)\n# but round-trips the wrong element name, so rebuild the paragraph's XML\n# directly (preserving its existing paragraph-mark attributes) and push it\n# back in with InsertXML.\n$para = $d.Paragraphs(1)\n$r = $para.Range\n\n$full = $d.Content.WordOpenXML\n$partMarker = 'pkg:name=\"/word/document.xml\"'\n$partIdx = $full.IndexOf($partMarker)\n$bodyIdx = $full.IndexOf(\"<w:body>\", $partIdx)\n$pStart = $full.IndexOf(\"<w:p\", $bodyIdx)\n$pEnd = $full.IndexOf(\">\", $pStart)\n$openTag = $full.Substring($pStart, $pEnd - $pStart + 1)\n\n$selfClosing = $openTag.EndsWith(\"/>\")\nif ($selfClosing) {\n  $attrs = $openTag.Substring(4, $openTag.Length - 4 - 2).Trim()\n} else {\n  $attrs = $openTag.Substring(4, $openTag.Length - 4 - 1).Trim()\n}\n\n$newParaXml = \"<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' $attrs>\" + `\n  \"<w:pPr><w:rPr><w:rStyle w:val='H5'/></w:rPr></w:pPr>\" + `\n  \"</w:p>\"\n\n$null = $r.InsertXML($newParaXml)\n\nWrite-Output \"done\"\n"}
